$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.978.27'
$ws.Range("E2").Value = '  -0.05%  '

$ws.Range("D3").Value = '2.272.94'
$ws.Range("E3").Value = '  +0.71%  '

$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '305.81'
$ws.Range("E5").Value = '  +1.19%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '93.03'
$ws.Range("E6").Value = '  +0.23%  '

$ws.Range("E7").Value = '  -0.29%  '

$ws.Range("E8").Value = '  +0.06%  '

$ws.Range("E9").Value = '  +1.02%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '32.78'
$ws.Range("E10").Value = '  +0.02%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0804'
$ws.Range("E11").Value = '  +0.27%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.113'
$ws.Range("E12").Value = '  -1.92%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.70'
$ws.Range("E13").Value = '  -0.04%  '

$ws.Range("D14").Value = '2.624.19'
$ws.Range("E14").Value = '  +0.74%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.37'
$ws.Range("E15").Value = '  +1.63%  '

$ws.Range("D16").Value = '2.278.85'
$ws.Range("E16").Value = '  +0.31%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.785'
$ws.Range("E17").Value = '  +3.60%  '

$ws.Range("D18").Value = '41.876.19'
$ws.Range("E18").Value = '  -0.02%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.79'
$ws.Range("E19").Value = '  +4.20%  '

$ws.Range("E20").Value = '  +1.61%  '

$ws.Range("E21").Value = '  +0.56%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.27'
$ws.Range("E22").Value = '  +1.58%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '244.40'
$ws.Range("E23").Value = '  +1.13%  '

$ws.Range("E24").Value = '  +0.22%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.94'
$ws.Range("E25").Value = '  +1.24%  '

$ws.Range("E26").Value = '  +0.03%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.01'
$ws.Range("E27").Value = '  +0.23%  '

$ws.Range("E28").Value = '  +0.12%  '

$ws.Range("E29").Value = '  -7.47%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.11'
$ws.Range("E30").Value = '  +2.66%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '159.59'
$ws.Range("E31").Value = '  +0.66%  '

$ws.Range("E32").Value = '  +3.86%  '

$ws.Range("E33").Value = '  +0.03%  '

$ws.Range("E34").Value = '  +0.13%  '

$ws.Range("E35").Value = '  -0.16%  '

$ws.Range("E36").Value = '  +3.52%  '

$ws.Range("E37").Value = '  -1.34%  '

$ws.Range("E38").Value = '  +0.16%  '

$ws.Range("E39").Value = '  +0.83%  '

$ws.Range("E40").Value = '  -0.11%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.96'
$ws.Range("E41").Value = '  +0.55%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '19.98'
$ws.Range("E42").Value = '  +0.34%  '

$ws.Range("D43").Value = '2.017.27'
$ws.Range("E43").Value = '  -1.65%  '

$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0283'
$ws.Range("E44").Value = '  +1.17%  '

$ws.Range("B45").Value = 'ApeXProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.22'
$ws.Range("E45").Value = '  +8.38%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.33'
$ws.Range("E46").Value = '  +1.81%  '

$ws.Range("E47").Value = '  +2.17%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '53.37'
$ws.Range("E48").Value = '  +2.78%  '

$ws.Range("E49").Value = '  -0.60%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '72.55'
$ws.Range("E50").Value = '  +2.97%  '

$ws.Range("E51").Value = '  +0.25%  '
